# Remove the empty paragraph, the long "Me gustaría añadir un tercer ..." paragraph,
# and the following empty paragraph that sit between "...lo contrario." and the next
# existing empty paragraph (which stays right before "Los espectros ...").
$d = $word.ActiveDocument

$startPara = $d.Paragraphs.Item(2)
$endPara = $d.Paragraphs.Item(4)

$r = $d.Range($startPara.Range.Start, $endPara.Range.End)
$r.Delete()
